# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (column G) values for rows 2-49, replacing the prior Strike# counts.
$kValues = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 1
    6  = 1
    7  = 1
    8  = 0
    9  = 1
    10 = 2
    11 = 0
    12 = 0
    13 = 0
    14 = 1
    15 = 2
    16 = 0
    17 = 1
    18 = 2
    19 = 1
    20 = 2
    21 = 2
    22 = 2
    23 = 0
    24 = 0
    25 = 0
    26 = 1
    27 = 0
    28 = 1
    29 = 1
    30 = 1
    31 = 1
    32 = 0
    33 = 1
    34 = 1
    35 = 1
    36 = 1
    37 = 0
    38 = 1
    39 = 1
    40 = 0
    41 = 0
    42 = 1
    43 = 0
    44 = 1
    45 = 1
    46 = 1
    47 = 1
    48 = 3
    49 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
